$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.87"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("D3").Value = "'41.71"
$ws.Range("E3").Value = "'4.22%"
$ws.Range("E4").Value = "'2.46%"
$ws.Range("D5").Value = "'0.07615"
$ws.Range("E5").Value = "'-0.79%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.252"
$ws.Range("E6").Value = "'0.25%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.626"
$ws.Range("E7").Value = "'0.17%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.497"
$ws.Range("E8").Value = "'-2.06%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9046"
$ws.Range("E9").Value = "'1.62%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1093"
$ws.Range("E10").Value = "'9.78%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1764"
$ws.Range("E11").Value = "'1.37%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09214"
$ws.Range("E12").Value = "'3.04%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04272"
$ws.Range("E13").Value = "'-2.47%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1051"
$ws.Range("E14").Value = "'-0.31%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001260"
$ws.Range("E15").Value = "'-0.97%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005822"
$ws.Range("E16").Value = "'-1.27%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.362"
$ws.Range("E17").Value = "'0.15%"
$ws.Range("E18").Value = "'-1.95%"
$ws.Range("D19").Value = "'6.554"
$ws.Range("E19").Value = "'-6.88%"
$ws.Range("D20").Value = "'0.1363"
$ws.Range("E20").Value = "'1.66%"
$ws.Range("E21").Value = "'-11.48%"
$ws.Range("E22").Value = "'-1.36%"
$ws.Range("D23").Value = "'0.001220"
$ws.Range("E23").Value = "'1.71%"
$ws.Range("D24").Value = "'0.004096"
$ws.Range("E24").Value = "'0.70%"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'6.50%"
$ws.Range("D38").Value = "'0.02421"
$ws.Range("E38").Value = "'2.59%"
$ws.Range("D39").Value = "'0.05197"
$ws.Range("E39").Value = "'0.44%"
$ws.Range("D40").Value = "'0.007775"
$ws.Range("E40").Value = "'-2.32%"
$ws.Range("D41").Value = "'0.1300"
$ws.Range("E41").Value = "'-1.72%"
$ws.Range("D42").Value = "'0.006955"
$ws.Range("E42").Value = "'5.88%"
$ws.Range("D43").Value = "'0.001921"
$ws.Range("E43").Value = "'-5.95%"
$ws.Range("D44").Value = "'0.008068"
$ws.Range("E44").Value = "'5.92%"
$ws.Range("D45").Value = "'0.3050"
$ws.Range("E45").Value = "'-0.26%"
$ws.Range("D46").Value = "'0.00006736"
$ws.Range("E46").Value = "'1.28%"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'0.009147"
$ws.Range("E48").Value = "'142.75%"
$ws.Range("E49").Value = "'-16.00%"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("E51").Value = "'-0.06%"
